$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.222.00"
$ws.Range("E2").Value = "  +2.70%  "

# Row 3
$ws.Range("D3").Value = "2.522.64"
$ws.Range("E3").Value = "  +1.79%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'322.68"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").Value = "'109.51"
$ws.Range("E6").Value = "  +1.76%  "

# Row 7
$ws.Range("D7").Value = "'0.534"
$ws.Range("E7").Value = "  +2.57%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.554"
$ws.Range("E9").Value = "  +4.20%  "

# Row 10
$ws.Range("D10").Value = "'40.48"
$ws.Range("E10").Value = "  +4.83%  "

# Row 11
$ws.Range("D11").Value = "'20.59"
$ws.Range("E11").Value = "  +13.89%  "

# Row 12
$ws.Range("D12").Value = "'0.0824"
$ws.Range("E12").Value = "  +2.32%  "

# Row 13
$ws.Range("E13").Value = "  +1.37%  "

# Row 14
$ws.Range("D14").Value = "'7.29"
$ws.Range("E14").Value = "  +2.66%  "

# Row 15
$ws.Range("D15").Value = "2.915.90"
$ws.Range("E15").Value = "  +1.78%  "

# Row 16
$ws.Range("D16").Value = "2.522.09"
$ws.Range("E16").Value = "  +2.04%  "

# Row 17
$ws.Range("D17").Value = "'0.856"
$ws.Range("E17").Value = "  +1.64%  "

# Row 18
$ws.Range("D18").Value = "48.060.75"
$ws.Range("E18").Value = "  +2.50%  "

# Row 19
$ws.Range("D19").Value = "'13.33"
$ws.Range("E19").Value = "  +5.02%  "

# Row 20
$ws.Range("D20").Value = "'6.63"
$ws.Range("E20").Value = "  +0.77%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0946"
$ws.Range("E21").Value = "  +1.95%  "

# Row 22
$ws.Range("D22").Value = "'2.71"
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
$ws.Range("D23").Value = "'72.02"
$ws.Range("E23").Value = "  +2.55%  "

# Row 24
$ws.Range("D24").Value = "'264.61"
$ws.Range("E24").Value = "  +8.23%  "

# Row 25
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +1.19%  "

# Row 26
$ws.Range("E26").Value = "  -0.25%  "

# Row 27
$ws.Range("D27").Value = "'26.12"
$ws.Range("E27").Value = "  +2.34%  "

# Row 28
$ws.Range("D28").Value = "'10.12"
$ws.Range("E28").Value = "  +1.09%  "

# Row 29
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.144"
$ws.Range("E29").Value = "  +4.43%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.20"
$ws.Range("E30").Value = "  -2.63%  "

# Row 31
$ws.Range("D31").Value = "'36.26"
$ws.Range("E31").Value = "  +4.30%  "

# Row 32
$ws.Range("D32").Value = "'49.70"
$ws.Range("E32").Value = "  +0.75%  "

# Row 33
$ws.Range("D33").Value = "'19.83"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").Value = "'5.42"
$ws.Range("E34").Value = "  +1.93%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").Value = "'0.0791"
$ws.Range("E36").Value = "  +1.79%  "

# Row 37
$ws.Range("D37").Value = "'1.99"
$ws.Range("E37").Value = "  +1.93%  "

# Row 38
$ws.Range("D38").Value = "'4.74"
$ws.Range("E38").Value = "  +2.44%  "

# Row 39
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +2.50%  "

# Row 40
$ws.Range("E40").Value = "  +0.89%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'120.38"
$ws.Range("E41").Value = "  +1.02%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'22.07"
$ws.Range("E42").Value = "  +1.42%  "

# Row 43
$ws.Range("E43").Value = "  -0.64%  "

# Row 44
$ws.Range("D44").Value = "'0.0301"
$ws.Range("E44").Value = "  +2.80%  "

# Row 45
$ws.Range("D45").Value = "2.015.55"
$ws.Range("E45").Value = "  +2.16%  "

# Row 46
$ws.Range("D46").Value = "'3.17"
$ws.Range("E46").Value = "  +5.99%  "

# Row 47
$ws.Range("D47").Value = "'1.92"
$ws.Range("E47").Value = "  +9.46%  "

# Row 48
$ws.Range("D48").Value = "'2.05"
$ws.Range("E48").Value = "  +2.77%  "

# Row 49
$ws.Range("D49").Value = "'9.13"
$ws.Range("E49").Value = "  +1.42%  "

# Row 50
$ws.Range("D50").Value = "'5.22"
$ws.Range("E50").Value = "  +2.35%  "

# Row 51
$ws.Range("D51").Value = "'78.99"
$ws.Range("E51").Value = "  +3.23%  "
